# Updated cryptos list (price / 1h volume change) refresh.
# For the "Price" column (D) many new values are plain numeric-looking
# strings (e.g. "308.96"); Excel's COM Value setter auto-detects those as
# numbers unless the cell is explicitly formatted as Text ("@") first, so
# we force text formatting only on the cells where that auto-detection
# would otherwise mangle the value (values that still contain two dots,
# like "26.922.00", are already unambiguous text and are left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.922.00"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.816.23"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.96"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4691"
$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("E8").Value = "  -1.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07374"
$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8716"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "1.831.33"
$ws.Range("E12").Value = "  +2.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.383"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07068"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.89"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008708"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "26.952.15"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.325"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  -1.91%  "

$ws.Range("D24").Value = "2.049.05"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.01"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.172"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.35"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.345"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.24"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7684"
$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.512"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.915"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.090"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05293"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.950"
$ws.Range("E40").Value = "  +1.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.275"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5323"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.356"
$ws.Range("E43").Value = "  -3.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1659"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.453"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4927"
$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.48"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.74"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06304"
$ws.Range("E51").Value = "  -0.43%  "
